$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status moved from "In Translation" to "Ready for handoff" (Overview summary
# columns for each locale, plus the per-locale Status column).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Refresh the generation / handoff timestamps to reflect the new report run.
$wsOverview.Range("G2").Value = "2016-08-23 02:38:51"
$wsZhCn.Range("H2").Value = "2016-08-23 02:38:46"

# The longer "Ready for handoff" status text needs more room, so widen the
# status columns to fit.
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
